# Fix the slide page-counter text boxes so they read "current/9"
# sequentially instead of the original duplicated "x/8" values.
$p = $ppt.ActivePresentation

$updates = @{
    2  = "1/9"
    3  = "2/9"
    4  = "3/9"
    5  = "4/9"
    6  = "5/9"
    7  = "6/9"
    8  = "7/9"
    9  = "8/9"
    10 = "9/9"
}

foreach ($slideIndex in $updates.Keys) {
    $newText = $updates[$slideIndex]
    $s = $p.Slides.Item($slideIndex)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t -match '^[0-9]+/[0-9]+$') {
                    $sh.TextFrame.TextRange.Text = $newText
                }
            }
        }
    }
}
